$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:B55")
$rng.Sort($ws.Range("B1"), 2, $null, $null, 1, $null, 1, 1, $false, $null, $null, 1)
